# Update the daily terminal gate pricing values.
# New day (Dec 17, 2025 -> serial 46008) prices are added at the top of each
# state/terminal block, existing rows shift down one day (Dec 16 -> previously
# top slot; Dec 13 data is dropped from the visible two-day window).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 46008
$ws.Range("D8").Value = 157.72
$ws.Range("E8").Value = 157.78
$ws.Range("F8").Value = 167.78
$ws.Range("G8").Value = 157.9
$ws.Range("A9").Value = 46008
$ws.Range("D9").Value = 157.72
$ws.Range("E9").Value = 157.78
$ws.Range("F9").Value = 167.78
$ws.Range("G9").Value = 157.9
$ws.Range("A10").Value = 46008
$ws.Range("D10").Value = 159.65
$ws.Range("E10").Value = 160.22
$ws.Range("F10").Value = 170.23
$ws.Range("G10").Value = 160.74
$ws.Range("A11").Value = 46007
$ws.Range("D11").Value = 159.01
$ws.Range("E11").Value = 158.24
$ws.Range("F11").Value = 168.24
$ws.Range("G11").Value = 158.36
$ws.Range("A12").Value = 46007
$ws.Range("D12").Value = 159.01
$ws.Range("E12").Value = 158.24
$ws.Range("F12").Value = 168.24
$ws.Range("G12").Value = 158.36
$ws.Range("A13").Value = 46007
$ws.Range("D13").Value = 160.69
$ws.Range("E13").Value = 160.87
$ws.Range("F13").Value = 170.87
$ws.Range("G13").Value = 161.39
$ws.Range("A17").Value = 46008
$ws.Range("D17").Value = 162.6
$ws.Range("E17").Value = 163.76
$ws.Range("F17").Value = 173.76
$ws.Range("A18").Value = 46007
$ws.Range("D18").Value = 163.64
$ws.Range("E18").Value = 164.17
$ws.Range("F18").Value = 174.17
$ws.Range("A22").Value = 46008
$ws.Range("D22").Value = 158.42
$ws.Range("E22").Value = 159.44
$ws.Range("F22").Value = 169.04
$ws.Range("G22").Value = 160.6
$ws.Range("A23").Value = 46008
$ws.Range("D23").Value = 164.45
$ws.Range("E23").Value = 164.55
$ws.Range("F23").Value = 174.55
$ws.Range("A24").Value = 46008
$ws.Range("D24").Value = 164.23
$ws.Range("E24").Value = 164.88
$ws.Range("F24").Value = 174.88
$ws.Range("A25").Value = 46008
$ws.Range("D25").Value = 164.73
$ws.Range("E25").Value = 164.3
$ws.Range("F25").Value = 174.3
$ws.Range("G25").Value = 164.07
$ws.Range("A26").Value = 46008
$ws.Range("D26").Value = 163.63
$ws.Range("E26").Value = 165.73
$ws.Range("F26").Value = 175.73
$ws.Range("A27").Value = 46007
$ws.Range("D27").Value = 159.49
$ws.Range("E27").Value = 159.98
$ws.Range("F27").Value = 169.57
$ws.Range("G27").Value = 161.14
$ws.Range("A28").Value = 46007
$ws.Range("D28").Value = 165.72
$ws.Range("E28").Value = 164.98
$ws.Range("F28").Value = 174.98
$ws.Range("A29").Value = 46007
$ws.Range("D29").Value = 165.49
$ws.Range("E29").Value = 165.3
$ws.Range("F29").Value = 175.3
$ws.Range("A30").Value = 46007
$ws.Range("D30").Value = 166.1
$ws.Range("E30").Value = 164.73
$ws.Range("F30").Value = 174.73
$ws.Range("G30").Value = 164.5
$ws.Range("A31").Value = 46007
$ws.Range("D31").Value = 164.9
$ws.Range("E31").Value = 166.16
$ws.Range("F31").Value = 176.16
$ws.Range("A35").Value = 46008
$ws.Range("D35").Value = 157.83
$ws.Range("E35").Value = 157.36
$ws.Range("F35").Value = 166.36
$ws.Range("A36").Value = 46007
$ws.Range("D36").Value = 158.87
$ws.Range("E36").Value = 157.79
$ws.Range("F36").Value = 166.79
$ws.Range("A40").Value = 46008
$ws.Range("D40").Value = 163.93
$ws.Range("E40").Value = 164.42
$ws.Range("F40").Value = 174.42
$ws.Range("A41").Value = 46008
$ws.Range("D41").Value = 163.65
$ws.Range("E41").Value = 164.84
$ws.Range("F41").Value = 174.84
$ws.Range("A42").Value = 46007
$ws.Range("D42").Value = 164.97
$ws.Range("E42").Value = 164.82
$ws.Range("F42").Value = 174.82
$ws.Range("A43").Value = 46007
$ws.Range("D43").Value = 164.69
$ws.Range("E43").Value = 165.24
$ws.Range("F43").Value = 175.24
$ws.Range("A47").Value = 46008
$ws.Range("D47").Value = 159.24
$ws.Range("E47").Value = 159.1
$ws.Range("F47").Value = 169.1
$ws.Range("A48").Value = 46008
$ws.Range("D48").Value = 159.05
$ws.Range("E48").Value = 159.19
$ws.Range("F48").Value = 169.19
$ws.Range("A49").Value = 46007
$ws.Range("D49").Value = 160.15
$ws.Range("E49").Value = 159.4
$ws.Range("F49").Value = 169.4
$ws.Range("A50").Value = 46007
$ws.Range("D50").Value = 159.96
$ws.Range("E50").Value = 159.49
$ws.Range("F50").Value = 169.49
$ws.Range("A54").Value = 46008
$ws.Range("D54").Value = 173.42
$ws.Range("E54").Value = 175.12
$ws.Range("F54").Value = 185.12
$ws.Range("A55").Value = 46008
$ws.Range("D55").Value = 161.62
$ws.Range("E55").Value = 170.12
$ws.Range("F55").Value = 180.12
$ws.Range("A56").Value = 46008
$ws.Range("D56").Value = 163.99
$ws.Range("A57").Value = 46008
$ws.Range("D57").Value = 163
$ws.Range("E57").Value = 164.39
$ws.Range("A58").Value = 46008
$ws.Range("D58").Value = 158.9
$ws.Range("E58").Value = 160.44
$ws.Range("F58").Value = 170.44
$ws.Range("A59").Value = 46008
$ws.Range("D59").Value = 165.51
$ws.Range("E59").Value = 172.37
$ws.Range("A60").Value = 46007
$ws.Range("D60").Value = 174.47
$ws.Range("E60").Value = 175.52
$ws.Range("F60").Value = 185.52
$ws.Range("A61").Value = 46007
$ws.Range("D61").Value = 162.66
$ws.Range("E61").Value = 170.66
$ws.Range("F61").Value = 180.66
$ws.Range("A62").Value = 46007
$ws.Range("D62").Value = 165.03
$ws.Range("A63").Value = 46007
$ws.Range("D63").Value = 164.04
$ws.Range("E63").Value = 164.93
$ws.Range("A64").Value = 46007
$ws.Range("D64").Value = 159.95
$ws.Range("E64").Value = 160.98
$ws.Range("F64").Value = 170.98
$ws.Range("A65").Value = 46007
$ws.Range("D65").Value = 166.56
$ws.Range("E65").Value = 172.78
